$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force text-typed (non-numeric-inferred) values,
# mirroring how the source workbook stores B-column org_id codes as text.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

# --- Update existing rows 2-31 (B and C columns) ---
# row 2
$scratch.Value2 = '4750'
$scratch.Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4163) | Out-Null
$ws.Range("C2").Value2 = 'Институт энергетической стратегии'

# row 3
$scratch.Value2 = '500'
$scratch.Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$ws.Range("C3").Value2 = 'Геологический институт РАН'

# row 4
$scratch.Value2 = '509'
$scratch.Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4163) | Out-Null
$ws.Range("C4").Value2 = 'АО «ТомскНИПИнефть»'

# row 5
$scratch.Value2 = '522'
$scratch.Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4163) | Out-Null
$ws.Range("C5").Value2 = 'АО “Российские космические системы”'

# row 6
$scratch.Value2 = '533'
$scratch.Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4163) | Out-Null
$ws.Range("C6").Value2 = 'ООО «Газпромнефть НТЦ»'

# row 7
$scratch.Value2 = '549'
$scratch.Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$ws.Range("C7").Value2 = 'Филиал ООО "ЛУКОЙЛ-Инжиниринг" "КогалымНИПИнефть" в г. Тюмени'

# row 8
$scratch.Value2 = '570'
$scratch.Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Value2 = 'ИПНГ РАН'

# row 9
$scratch.Value2 = '5932'
$scratch.Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4163) | Out-Null
$ws.Range("C9").Value2 = 'I.M. Sechenov First Moscow State Medical University, the Ministry of Health of the Russian Federation (Sechenov University)'

# row 10
$scratch.Value2 = '6019'
$scratch.Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4163) | Out-Null
$ws.Range("C10").Value2 = 'Институт теории прогноза землетрясений и математической геофизики РАН'

# row 11
$scratch.Value2 = '6228'
$scratch.Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4163) | Out-Null
$ws.Range("C11").Value2 = 'Almaty Institute of Power Engineering and Telecommunications'

# row 12
$scratch.Value2 = '6317'
$scratch.Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4163) | Out-Null
$ws.Range("C12").Value2 = 'НТС ПАО "Газпром"'

# row 13
$scratch.Value2 = '6939'
$scratch.Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null
$ws.Range("C13").Value2 = 'ПАО «Татнефть»'

# row 14
$scratch.Value2 = '7109'
$scratch.Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4163) | Out-Null
$ws.Range("C14").Value2 = 'ПАО «НК «Роснефть»'

# row 15
$scratch.Value2 = '734'
$scratch.Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("C15").Value2 = 'Ammosov North Eastern Federal University'

# row 16
$scratch.Value2 = '7655'
$scratch.Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4163) | Out-Null
$ws.Range("C16").Value2 = 'Московский государственный университет им. М.В. Ломоносова'

# row 17
$scratch.Value2 = '80'
$scratch.Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4163) | Out-Null
$ws.Range("C17").Value2 = 'V.A. Trapeznikov Institute of Control Sciences|Russian Academy of Sciences'

# row 18
$scratch.Value2 = '848'
$scratch.Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4163) | Out-Null
$ws.Range("C18").Value2 = 'Institute of Geography RAS'

# row 19
$scratch.Value2 = '870'
$scratch.Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4163) | Out-Null
$ws.Range("C19").Value2 = 'The Moscow Mining University, NUSTMIS&S'

# row 20
$scratch.Value2 = '898'
$scratch.Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4163) | Out-Null
$ws.Range("C20").Value2 = 'Российский химико-технологический университет им. Д.И. Менделеева'

# row 21
$scratch.Value2 = '982'
$scratch.Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4163) | Out-Null
$ws.Range("C21").Value2 = 'Институт микробиологии им. С.Н. Виноградского, ФИЦ Биотехнологии РАН'

# row 22
$scratch.Value2 = '4765'
$scratch.Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4163) | Out-Null
$ws.Range("C22").Value2 = 'Институт морской геологии и геофизики ДВО РАН'

# row 23
$scratch.Value2 = '452'
$scratch.Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4163) | Out-Null
$ws.Range("C23").Value2 = 'Институт нефтехимического синтеза им. А.В. Топчиева РАН'

# row 24
$scratch.Value2 = '11374'
$scratch.Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4163) | Out-Null
$ws.Range("C24").Value2 = 'Ухтинский государственный технический университет'

# row 25
$scratch.Value2 = '440'
$scratch.Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4163) | Out-Null
$ws.Range("C25").Value2 = 'Norilsk State Industrial Institute'

# row 26
$scratch.Value2 = '132'
$scratch.Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4163) | Out-Null
$ws.Range("C26").Value2 = 'ФГАОУ ВО «Российский университет дружбы народов»'

# row 27
$scratch.Value2 = '13814'
$scratch.Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4163) | Out-Null
$ws.Range("C27").Value2 = 'A.N. Nesmeyanov Institute of Organoelement Compounds of Russian Academy of Sciences'

# row 28
$scratch.Value2 = '14430'
$scratch.Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4163) | Out-Null
$ws.Range("C28").Value2 = 'Prokhorov General Physics Institute of the Russian Academy of Sciences'

# row 29
$scratch.Value2 = '14552'
$scratch.Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4163) | Out-Null
$ws.Range("C29").Value2 = 'Northern (Arctic) Federal University named after M.V. Lomonosov'

# row 30
$scratch.Value2 = '14819'
$scratch.Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4163) | Out-Null
$ws.Range("C30").Value2 = 'Тюменский государственный университет'

# row 31
$scratch.Value2 = '15203'
$scratch.Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4163) | Out-Null
$ws.Range("C31").Value2 = 'Институт проблем нефти и газа РАН'

# --- Append new rows 32-70 (A, B, C columns) ---
$aFormatSource = $ws.Range("A2")
# row 32
$aFormatSource.Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value2 = 30
$scratch.Value2 = '162'
$scratch.Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4163) | Out-Null
$ws.Range("C32").Value2 = 'Всероссийский научно-исследовательский геологический институт им. А.П. Карпинского'

# row 33
$aFormatSource.Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Value2 = 31
$scratch.Value2 = '176'
$scratch.Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4163) | Out-Null
$ws.Range("C33").Value2 = 'ООО «Газпром добыча Надым»'

# row 34
$aFormatSource.Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("A34").Value2 = 32
$scratch.Value2 = '19244'
$scratch.Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4163) | Out-Null
$ws.Range("C34").Value2 = 'ООО «РН-СахалинНИПИморнефть»'

# row 35
$aFormatSource.Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Range("A35").Value2 = 33
$scratch.Value2 = '2258'
$scratch.Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4163) | Out-Null
$ws.Range("C35").Value2 = 'Институт криосферы Земли ТюмНЦ СО РАН'

# row 36
$aFormatSource.Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").Value2 = 34
$scratch.Value2 = '2268'
$scratch.Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4163) | Out-Null
$ws.Range("C36").Value2 = 'ВНИИОкеангеология'

# row 37
$aFormatSource.Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null
$ws.Range("A37").Value2 = 35
$scratch.Value2 = '2280'
$scratch.Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4163) | Out-Null
$ws.Range("C37").Value2 = 'All-Russian Research Geological Oil Institute'

# row 38
$aFormatSource.Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value2 = 36
$scratch.Value2 = '252'
$scratch.Copy() | Out-Null
$ws.Range("B38").PasteSpecial(-4163) | Out-Null
$ws.Range("C38").Value2 = 'V.B. Sochava Institute of Geography|Siberian Branch of Russian Academy of Sciences'

# row 39
$aFormatSource.Copy() | Out-Null
$ws.Range("A39").PasteSpecial(-4122) | Out-Null
$ws.Range("A39").Value2 = 37
$scratch.Value2 = '2541'
$scratch.Copy() | Out-Null
$ws.Range("B39").PasteSpecial(-4163) | Out-Null
$ws.Range("C39").Value2 = 'Альметьевский гос. нефтяной институт'

# row 40
$aFormatSource.Copy() | Out-Null
$ws.Range("A40").PasteSpecial(-4122) | Out-Null
$ws.Range("A40").Value2 = 38
$scratch.Value2 = '257'
$scratch.Copy() | Out-Null
$ws.Range("B40").PasteSpecial(-4163) | Out-Null
$ws.Range("C40").Value2 = 'Frumkin Institute of Physical Chemistry and Electrochemistry|Russian Academy of Sciences'

# row 41
$aFormatSource.Copy() | Out-Null
$ws.Range("A41").PasteSpecial(-4122) | Out-Null
$ws.Range("A41").Value2 = 39
$scratch.Value2 = '327'
$scratch.Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4163) | Out-Null
$ws.Range("C41").Value2 = 'Пермский национальный исследовательский политехнический университет'

# row 42
$aFormatSource.Copy() | Out-Null
$ws.Range("A42").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Value2 = 40
$scratch.Value2 = '361'
$scratch.Copy() | Out-Null
$ws.Range("B42").PasteSpecial(-4163) | Out-Null
$ws.Range("C42").Value2 = 'Landau Institute for Theoretical Physics|Russian Academy of Sciences'

# row 43
$aFormatSource.Copy() | Out-Null
$ws.Range("A43").PasteSpecial(-4122) | Out-Null
$ws.Range("A43").Value2 = 41
$scratch.Value2 = '364'
$scratch.Copy() | Out-Null
$ws.Range("B43").PasteSpecial(-4163) | Out-Null
$ws.Range("C43").Value2 = 'Gubkin University'

# row 44
$aFormatSource.Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null
$ws.Range("A44").Value2 = 42
$scratch.Value2 = '371'
$scratch.Copy() | Out-Null
$ws.Range("B44").PasteSpecial(-4163) | Out-Null
$ws.Range("C44").Value2 = 'Российский государственный геологоразведочный университет имени Серго Орджоникидзе'

# row 45
$aFormatSource.Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Range("A45").Value2 = 43
$scratch.Value2 = '99'
$scratch.Copy() | Out-Null
$ws.Range("B45").PasteSpecial(-4163) | Out-Null
$ws.Range("C45").Value2 = 'Институт физики Земли им. О.Ю. Шмидта РАН'

# row 46
$aFormatSource.Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null
$ws.Range("A46").Value2 = 44
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B46").PasteSpecial(-4163) | Out-Null
$ws.Range("C46").Value2 = '18 Vinogradnaya Street, Alushta, 298517'

# row 47
$aFormatSource.Copy() | Out-Null
$ws.Range("A47").PasteSpecial(-4122) | Out-Null
$ws.Range("A47").Value2 = 45
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B47").PasteSpecial(-4163) | Out-Null
$ws.Range("C47").Value2 = 'Weatherford'

# row 48
$aFormatSource.Copy() | Out-Null
$ws.Range("A48").PasteSpecial(-4122) | Out-Null
$ws.Range("A48").Value2 = 46
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B48").PasteSpecial(-4163) | Out-Null
$ws.Range("C48").Value2 = 'Institute of Energy Strategy'

# row 49
$aFormatSource.Copy() | Out-Null
$ws.Range("A49").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Value2 = 47
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B49").PasteSpecial(-4163) | Out-Null
$ws.Range("C49").Value2 = 'ООО НТФ "Атомбиотех"'

# row 50
$aFormatSource.Copy() | Out-Null
$ws.Range("A50").PasteSpecial(-4122) | Out-Null
$ws.Range("A50").Value2 = 48
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B50").PasteSpecial(-4163) | Out-Null
$ws.Range("C50").Value2 = 'Oil and Gas Research Institute Russian Academy of Sciences (IPNG RAS)'

# row 51
$aFormatSource.Copy() | Out-Null
$ws.Range("A51").PasteSpecial(-4122) | Out-Null
$ws.Range("A51").Value2 = 49
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B51").PasteSpecial(-4163) | Out-Null
$ws.Range("C51").Value2 = 'IPNG RAN'

# row 52
$aFormatSource.Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null
$ws.Range("A52").Value2 = 50
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B52").PasteSpecial(-4163) | Out-Null
$ws.Range("C52").Value2 = 'Institute of Oil and Gas Problems of the Russian Academy of Sciences'

# row 53
$aFormatSource.Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null
$ws.Range("A53").Value2 = 51
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B53").PasteSpecial(-4163) | Out-Null
$ws.Range("C53").Value2 = 'АО "Томский научно-исследовательский и проектный институт нефти и газа" (АО "ТомскНИПИнефть")'

# row 54
$aFormatSource.Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$ws.Range("A54").Value2 = 52
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B54").PasteSpecial(-4163) | Out-Null
$ws.Range("C54").Value2 = 'Institute of Solid State Physics|Russian Academy of Sciences Academician'

# row 55
$aFormatSource.Copy() | Out-Null
$ws.Range("A55").PasteSpecial(-4122) | Out-Null
$ws.Range("A55").Value2 = 53
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B55").PasteSpecial(-4163) | Out-Null
$ws.Range("C55").Value2 = 'ООО «Тюменский институт нефти и газа»'

# row 56
$aFormatSource.Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null
$ws.Range("A56").Value2 = 54
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B56").PasteSpecial(-4163) | Out-Null
$ws.Range("C56").Value2 = 'Группа компаний ITPS'

# row 57
$aFormatSource.Copy() | Out-Null
$ws.Range("A57").PasteSpecial(-4122) | Out-Null
$ws.Range("A57").Value2 = 55
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B57").PasteSpecial(-4163) | Out-Null
$ws.Range("C57").Value2 = 'Wildcat Technologies LLC'

# row 58
$aFormatSource.Copy() | Out-Null
$ws.Range("A58").PasteSpecial(-4122) | Out-Null
$ws.Range("A58").Value2 = 56
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B58").PasteSpecial(-4163) | Out-Null
$ws.Range("C58").Value2 = 'Хромос Инжиниринг'

# row 59
$aFormatSource.Copy() | Out-Null
$ws.Range("A59").PasteSpecial(-4122) | Out-Null
$ws.Range("A59").Value2 = 57
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B59").PasteSpecial(-4163) | Out-Null
$ws.Range("C59").Value2 = 'Department of Biotechnology|I. M. Sechenov First Moscow State Medical University'

# row 60
$aFormatSource.Copy() | Out-Null
$ws.Range("A60").PasteSpecial(-4122) | Out-Null
$ws.Range("A60").Value2 = 58
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B60").PasteSpecial(-4163) | Out-Null
$ws.Range("C60").Value2 = 'Landau Institute for Theoretical Physics of the RAS'

# row 61
$aFormatSource.Copy() | Out-Null
$ws.Range("A61").PasteSpecial(-4122) | Out-Null
$ws.Range("A61").Value2 = 59
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B61").PasteSpecial(-4163) | Out-Null
$ws.Range("C61").Value2 = 'Institute of Solid State Physics of the RAS'

# row 62
$aFormatSource.Copy() | Out-Null
$ws.Range("A62").PasteSpecial(-4122) | Out-Null
$ws.Range("A62").Value2 = 60
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B62").PasteSpecial(-4163) | Out-Null
$ws.Range("C62").Value2 = 'ООО "Хромос Инжиниринг"'

# row 63
$aFormatSource.Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4122) | Out-Null
$ws.Range("A63").Value2 = 61
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B63").PasteSpecial(-4163) | Out-Null
$ws.Range("C63").Value2 = 'Gazprom Export'

# row 64
$aFormatSource.Copy() | Out-Null
$ws.Range("A64").PasteSpecial(-4122) | Out-Null
$ws.Range("A64").Value2 = 62
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B64").PasteSpecial(-4163) | Out-Null
$ws.Range("C64").Value2 = 'Scientific Council of RAS on System Research in Energy'

# row 65
$aFormatSource.Copy() | Out-Null
$ws.Range("A65").PasteSpecial(-4122) | Out-Null
$ws.Range("A65").Value2 = 63
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B65").PasteSpecial(-4163) | Out-Null
$ws.Range("C65").Value2 = 'Nansen Environmental and Remote Sensing Centre|Bjerknes Centre for Climate Research'

# row 66
$aFormatSource.Copy() | Out-Null
$ws.Range("A66").PasteSpecial(-4122) | Out-Null
$ws.Range("A66").Value2 = 64
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B66").PasteSpecial(-4163) | Out-Null
$ws.Range("C66").Value2 = 'Department of Geography|The George Washington University'

# row 67
$aFormatSource.Copy() | Out-Null
$ws.Range("A67").PasteSpecial(-4122) | Out-Null
$ws.Range("A67").Value2 = 65
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B67").PasteSpecial(-4163) | Out-Null
$ws.Range("C67").Value2 = 'Межотраслевой экспертно - аналитический центр Союза Нефтегазопромышленников России'

# row 68
$aFormatSource.Copy() | Out-Null
$ws.Range("A68").PasteSpecial(-4122) | Out-Null
$ws.Range("A68").Value2 = 66
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B68").PasteSpecial(-4163) | Out-Null
$ws.Range("C68").Value2 = 'ООО «Красноярскгазпром нефтегазпроект»'

# row 69
$aFormatSource.Copy() | Out-Null
$ws.Range("A69").PasteSpecial(-4122) | Out-Null
$ws.Range("A69").Value2 = 67
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B69").PasteSpecial(-4163) | Out-Null
$ws.Range("C69").Value2 = 'ООО «КБ Стрелка»'

# row 70
$aFormatSource.Copy() | Out-Null
$ws.Range("A70").PasteSpecial(-4122) | Out-Null
$ws.Range("A70").Value2 = 68
$scratch.Value2 = ' '
$scratch.Copy() | Out-Null
$ws.Range("B70").PasteSpecial(-4163) | Out-Null
$ws.Range("C70").Value2 = 'ООО Тимано-Печорский научно-исследовательский центр (ООО «ТПНИЦ»)'

# Clean up the scratch cell
$scratch.Clear() | Out-Null
$excel.CutCopyMode = 0